# Decrease column E ("剩余") by 1 for every data row (2-99) except row 36,
# which stays unchanged per the source diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 99; $row++) {
    if ($row -eq 36) {
        continue
    }
    $cell = $ws.Cells.Item($row, 5)   # column E = 5
    $current = $cell.Value2
    if ($null -ne $current) {
        $cell.Value2 = $current - 1
    }
}
